$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.249.42"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.227.57"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.58%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.43"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -6.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.580"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -7.58%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.559"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -7.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.03"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.03"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("E12").Value = "  -8.91%  "
$ws.Range("E13").Value = "  -7.07%  "
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.566.53"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.862"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -10.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.33"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.232.47"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.156.08"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -8.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0967"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -8.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.55"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -9.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.21"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -11.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.15"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -10.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "236.71"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -7.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.18"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.03%  "
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.05"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.03"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -9.71%  "
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.41"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -12.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "36.70"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.24"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -7.36%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.30"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.06%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0865"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.67"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.26"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.68%  "
$ws.Range("E38").Value = "  -7.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.41"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.83"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.44%  "
$ws.Range("E41").Value = "  -8.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.71"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0319"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -8.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.36"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +13.06%  "
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.772.31"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.202"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -9.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "83.55"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -11.81%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.86"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.38%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.27"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -11.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.50"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -11.52%  "
